$d = $word.ActiveDocument

# 1) Remove the "Meta description: ..." paragraph that currently sits right
#    after the title ("Play Big Bad Wolf Free Slot Game - Review and Bonuses").
$titlePara = $d.Paragraphs.Item(1)
$metaPara = $titlePara.Next()
$metaPara.Range.Delete()

# 2) Before the final paragraph (currently the DALLE image prompt text),
#    insert a new paragraph containing the bolded title text
#    "Play Big Bad Wolf Free Slot Game - Review and Bonuses". Inserting it
#    *after* the preceding (non-italic) paragraph -- rather than *before*
#    the italic prompt paragraph -- keeps it from inheriting italics.
$count = $d.Paragraphs.Count
$precedingPara = $d.Paragraphs.Item($count - 1)
$precedingPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($count)
$newPara.Style = "Normal"
$newPara.Range.Text = "Play Big Bad Wolf Free Slot Game - Review and Bonuses"

$start = $newPara.Range.Start
$end = $newPara.Range.End - 1
$textRange = $d.Range($start, $end)
$textRange.Font.Bold = $true

# 3) Replace the DALLE prompt text (now the last paragraph) with the
#    meta-description copy, keeping its existing italic run formatting.
$oldText = 'Prompt: Create a feature image fitting the game "Big Bad Wolf". DALLE, please create a cartoon-style feature image for the game "Big Bad Wolf" that features a happy Maya warrior with glasses. The Maya warrior should be holding a basket of apples and standing in front of a countryside landscape with hills and a straw house in the background. The image should also include the Wolf and the Three Little Pigs as cartoon characters. The setting should be under the moonlight, with stars shining brightly in the sky. The colors should be fun and vibrant, with a focus on shades of blue and yellow. Can''t wait to see your creative work!'
$newText = "Immerse yourself in the beautiful countryside and win big with Big Bad Wolf slot game. Read our review and play for free now. Bonuses included."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
